# Update the NATMI ligand-receptor (Cntn2-Cntn1) results sheet with refreshed
# TPM-based values, and add the new Inflammatory-Mac sending/target cluster
# combinations (rows 4-7) alongside the existing FAPs/MuSCs rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cntn2"
$ws.Range("C2").Value = "Cntn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1436673333333333
$ws.Range("H2").Value = 0.431002
$ws.Range("I2").Value = 0.5549265272962071
$ws.Range("J2").Value = 0.5549265272962071
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1199133333333333
$ws.Range("N2").Value = 0.35974
$ws.Range("O2").Value = 0.3099390012751145
$ws.Range("P2").Value = 0.3099390012751145
$ws.Range("Q2").Value = 0.01722762883111111
$ws.Range("R2").Value = 0.15504865948
$ws.Range("S2").Value = 0.171993373651254
$ws.Range("T2").Value = 0.171993373651254

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cntn2"
$ws.Range("C3").Value = "Cntn1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1436673333333333
$ws.Range("H3").Value = 0.431002
$ws.Range("I3").Value = 0.5549265272962071
$ws.Range("J3").Value = 0.5549265272962071
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2669800000000001
$ws.Range("N3").Value = 0.8009400000000001
$ws.Range("O3").Value = 0.6900609987248855
$ws.Range("P3").Value = 0.6900609987248854
$ws.Range("Q3").Value = 0.03835630465333334
$ws.Range("R3").Value = 0.3452067418800001
$ws.Range("S3").Value = 0.3829331536449531
$ws.Range("T3").Value = 0.382933153644953

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Cntn2"
$ws.Range("C4").Value = "Cntn1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.011717
$ws.Range("H4").Value = 0.035151
$ws.Range("I4").Value = 0.04525784650880733
$ws.Range("J4").Value = 0.04525784650880733
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1199133333333333
$ws.Range("N4").Value = 0.35974
$ws.Range("O4").Value = 0.3099390012751145
$ws.Range("P4").Value = 0.3099390012751145
$ws.Range("Q4").Value = 0.001405024526666667
$ws.Range("R4").Value = 0.01264522074
$ws.Range("S4").Value = 0.01402717174680217
$ws.Range("T4").Value = 0.01402717174680217

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Cntn2"
$ws.Range("C5").Value = "Cntn1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.011717
$ws.Range("H5").Value = 0.035151
$ws.Range("I5").Value = 0.04525784650880733
$ws.Range("J5").Value = 0.04525784650880733
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2669800000000001
$ws.Range("N5").Value = 0.8009400000000001
$ws.Range("O5").Value = 0.6900609987248855
$ws.Range("P5").Value = 0.6900609987248854
$ws.Range("Q5").Value = 0.003128204660000001
$ws.Range("R5").Value = 0.02815384194
$ws.Range("S5").Value = 0.03123067476200515
$ws.Range("T5").Value = 0.03123067476200515

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cntn2"
$ws.Range("C6").Value = "Cntn1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.10351
$ws.Range("H6").Value = 0.31053
$ws.Range("I6").Value = 0.3998156261949856
$ws.Range("J6").Value = 0.3998156261949856
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1199133333333333
$ws.Range("N6").Value = 0.35974
$ws.Range("O6").Value = 0.3099390012751145
$ws.Range("P6").Value = 0.3099390012751145
$ws.Range("Q6").Value = 0.01241222913333333
$ws.Range("R6").Value = 0.1117100622
$ws.Range("S6").Value = 0.1239184558770583
$ws.Range("T6").Value = 0.1239184558770583

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cntn2"
$ws.Range("C7").Value = "Cntn1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.10351
$ws.Range("H7").Value = 0.31053
$ws.Range("I7").Value = 0.3998156261949856
$ws.Range("J7").Value = 0.3998156261949856
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2669800000000001
$ws.Range("N7").Value = 0.8009400000000001
$ws.Range("O7").Value = 0.6900609987248855
$ws.Range("P7").Value = 0.6900609987248854
$ws.Range("Q7").Value = 0.0276350998
$ws.Range("R7").Value = 0.2487158982
$ws.Range("S7").Value = 0.2758971703179272
$ws.Range("T7").Value = 0.2758971703179272
